$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '43.238.80'
Set-TextCell $ws.Range("E2") '  +0.29%  '

# Row 3
Set-TextCell $ws.Range("D3") '2.362.09'
Set-TextCell $ws.Range("E3") '  +2.03%  '

# Row 4
Set-TextCell $ws.Range("E4") '  +0.03%  '

# Row 5
Set-TextCell $ws.Range("D5") '309.28'
Set-TextCell $ws.Range("E5") '  -0.27%  '

# Row 6
Set-TextCell $ws.Range("D6") '103.46'
Set-TextCell $ws.Range("E6") '  +2.42%  '

# Row 7
Set-TextCell $ws.Range("E7") '  -4.53%  '

# Row 8
Set-TextCell $ws.Range("E8") '  +0.02%  '

# Row 9
Set-TextCell $ws.Range("D9") '0.521'
Set-TextCell $ws.Range("E9") '  -0.88%  '

# Row 10
Set-TextCell $ws.Range("D10") '35.59'
Set-TextCell $ws.Range("E10") '  -1.26%  '

# Row 11
Set-TextCell $ws.Range("D11") '52.89'
Set-TextCell $ws.Range("E11") '  +1.35%  '

# Row 12
Set-TextCell $ws.Range("D12") '0.0804'
Set-TextCell $ws.Range("E12") '  -1.42%  '

# Row 13
Set-TextCell $ws.Range("E13") '  -0.47%  '

# Row 14
Set-TextCell $ws.Range("D14") '6.90'
Set-TextCell $ws.Range("E14") '  -4.07%  '

# Row 15
Set-TextCell $ws.Range("D15") '2.734.33'
Set-TextCell $ws.Range("E15") '  +2.40%  '

# Row 16
Set-TextCell $ws.Range("D16") '15.52'
Set-TextCell $ws.Range("E16") '  +3.52%  '

# Row 17
Set-TextCell $ws.Range("D17") '2.363.98'
Set-TextCell $ws.Range("E17") '  +2.29%  '

# Row 18
Set-TextCell $ws.Range("E18") '  -0.76%  '

# Row 19
Set-TextCell $ws.Range("D19") '43.221.99'
Set-TextCell $ws.Range("E19") '  +0.45%  '

# Row 20
Set-TextCell $ws.Range("D20") '6.30'
Set-TextCell $ws.Range("E20") '  +3.01%  '

# Row 21
Set-TextCell $ws.Range("D21") '11.84'
Set-TextCell $ws.Range("E21") '  -5.76%  '

# Row 22
Set-TextCell $ws.Range("E22") '  -0.91%  '

# Row 23
Set-TextCell $ws.Range("D23") '68.01'
Set-TextCell $ws.Range("E23") '  -0.64%  '

# Row 24
Set-TextCell $ws.Range("D24") '239.38'
Set-TextCell $ws.Range("E24") '  -0.74%  '

# Row 25
Set-TextCell $ws.Range("E25") '  +0.24%  '

# Row 26
Set-TextCell $ws.Range("E26") '  -1.49%  '

# Row 27
Set-TextCell $ws.Range("E27") '  +0.02%  '

# Row 28
Set-TextCell $ws.Range("D28") '25.62'
Set-TextCell $ws.Range("E28") '  +3.74%  '

# Row 29
Set-TextCell $ws.Range("B29") 'Toncoin'
Set-TextCell $ws.Range("C29") 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws.Range("D29") '2.33'
Set-TextCell $ws.Range("E29") '  +9.95%  '

# Row 30
Set-TextCell $ws.Range("B30") 'InjectiveProtocol'
Set-TextCell $ws.Range("C30") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws.Range("D30") '36.36'
Set-TextCell $ws.Range("E30") '  -2.93%  '

# Row 31
Set-TextCell $ws.Range("B31") 'Cosmos'
Set-TextCell $ws.Range("C31") 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell $ws.Range("D31") '9.44'
Set-TextCell $ws.Range("E31") '  -2.42%  '

# Row 32
Set-TextCell $ws.Range("B32") 'Monero'
Set-TextCell $ws.Range("C32") 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextCell $ws.Range("D32") '161.37'
Set-TextCell $ws.Range("E32") '  -3.31%  '

# Row 33
Set-TextCell $ws.Range("B33") 'Filecoin'
Set-TextCell $ws.Range("C33") 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell $ws.Range("D33") '5.21'
Set-TextCell $ws.Range("E33") '  -2.37%  '

# Row 34
Set-TextCell $ws.Range("B34") 'FirstDigitalUSD'
Set-TextCell $ws.Range("C34") 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextCell $ws.Range("D34") '1.00'
Set-TextCell $ws.Range("E34") '  -0.01%  '

# Row 35
Set-TextCell $ws.Range("B35") 'Celestia'
Set-TextCell $ws.Range("C35") 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextCell $ws.Range("D35") '18.11'
Set-TextCell $ws.Range("E35") '  +0.33%  '

# Row 36
Set-TextCell $ws.Range("B36") 'WEMIXToken'
Set-TextCell $ws.Range("C36") 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextCell $ws.Range("D36") '2.49'
Set-TextCell $ws.Range("E36") '  +4.17%  '

# Row 37
Set-TextCell $ws.Range("B37") 'RenderToken'
Set-TextCell $ws.Range("C37") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell $ws.Range("D37") '4.66'
Set-TextCell $ws.Range("E37") '  +8.53%  '

# Row 38
Set-TextCell $ws.Range("D38") '3.06'
Set-TextCell $ws.Range("E38") '  -3.61%  '

# Row 39
Set-TextCell $ws.Range("B39") 'Hedera'
Set-TextCell $ws.Range("C39") 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell $ws.Range("D39") '0.0735'
Set-TextCell $ws.Range("E39") '  -1.25%  '

# Row 40
Set-TextCell $ws.Range("B40") 'ARBITRUM'
Set-TextCell $ws.Range("C40") 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextCell $ws.Range("D40") '1.91'
Set-TextCell $ws.Range("E40") '  +3.24%  '

# Row 41
Set-TextCell $ws.Range("B41") 'Kaspa'
Set-TextCell $ws.Range("C41") 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell $ws.Range("D41") '0.104'
Set-TextCell $ws.Range("E41") '  -2.46%  '

# Row 42
Set-TextCell $ws.Range("B42") 'Stellar'
Set-TextCell $ws.Range("C42") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws.Range("D42") '0.113'
Set-TextCell $ws.Range("E42") '  -2.53%  '

# Row 43
Set-TextCell $ws.Range("B43") 'ApeXProtocol'
Set-TextCell $ws.Range("C43") 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell $ws.Range("D43") '2.59'
Set-TextCell $ws.Range("E43") '  +11.99%  '

# Row 44
Set-TextCell $ws.Range("B44") 'Maker'
Set-TextCell $ws.Range("C44") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws.Range("D44") '2.038.13'
Set-TextCell $ws.Range("E44") '  +3.16%  '

# Row 45
Set-TextCell $ws.Range("B45") 'EnergySwap'
Set-TextCell $ws.Range("C45") 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell $ws.Range("D45") '19.56'
Set-TextCell $ws.Range("E45") '  -0.08%  '

# Row 46
Set-TextCell $ws.Range("B46") 'VeChain'
Set-TextCell $ws.Range("C46") 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell $ws.Range("D46") '0.0288'
Set-TextCell $ws.Range("E46") '  -0.70%  '

# Row 47
Set-TextCell $ws.Range("B47") 'FraxShare'
Set-TextCell $ws.Range("C47") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws.Range("D47") '10.58'
Set-TextCell $ws.Range("E47") '  +7.87%  '

# Row 48
Set-TextCell $ws.Range("B48") 'NEARProtocol'
Set-TextCell $ws.Range("C48") 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextCell $ws.Range("D48") '3.07'
Set-TextCell $ws.Range("E48") '  +1.77%  '

# Row 49
Set-TextCell $ws.Range("B49") 'MultiversX'
Set-TextCell $ws.Range("C49") 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
Set-TextCell $ws.Range("D49") '57.54'
Set-TextCell $ws.Range("E49") '  +3.19%  '

# Row 50
Set-TextCell $ws.Range("B50") 'HuobiToken'
Set-TextCell $ws.Range("C50") 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell $ws.Range("D50") '2.92'
Set-TextCell $ws.Range("E50") '  -1.94%  '

# Row 51
Set-TextCell $ws.Range("B51") 'RocketPoolETH'
Set-TextCell $ws.Range("C51") 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextCell $ws.Range("D51") '2.595.48'
Set-TextCell $ws.Range("E51") '  +2.25%  '
